# Update Sheet2 ("browser"/"chrome" quick-reference table) to the new
# action/value reference table used by the automation suite.
#
# Target layout (Sheet2):
#   A1: action    B1: prompt
#   A2: accept    B2: Hello people!!
#   A3: dismiss   B3: hellll
#
# Values are written in this specific order (A1, A2, B1, B2, A3, B3) so
# that new shared-string entries are appended in the same sequence as the
# authoritative edit (action, accept, prompt, Hello people!!, dismiss,
# hellll).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1").Value = "action"
$ws2.Range("A2").Value = "accept"
$ws2.Range("B1").Value = "prompt"
$ws2.Range("B2").Value = "Hello people!!"
$ws2.Range("A3").Value = "dismiss"
$ws2.Range("B3").Value = "hellll"

# Column B needs to be widened to fit "Hello people!!" (best-fit column).
$ws2.Columns.Item(2).ColumnWidth = 11.6

# Sheet2 becomes the active/selected sheet, with M6 as the selected cell
# (Sheet1 loses the tabSelected flag automatically when Sheet2 activates).
$ws2.Activate()
$ws2.Range("M6").Select()
